$d = $word.ActiveDocument

# Each header-row cell below gets new text AND is turned Bold (w:b w:val="0" -> w:b).
# Using Find to locate the run's range, replacing its Text, then setting Font.Bold
# keeps the rest of the run's existing character formatting untouched.

$replacements = @(
    @{ Old = "차이 티 총 판매량(개)";         New = "총 차이 판매액(단위)" },
    @{ Old = "수제 차이 티 판매량(개)";        New = "Artisanal Chai 판매(단위)" },
    @{ Old = "즉석 음용 차이 티 판매량(개)";    New = "미리 만든 Chai 판매(단위)" },
    @{ Old = "소셜 미디어 참여율(조회수)";      New = "소셜 미디어 참여도(보기)" },
    @{ Old = "온라인 차이 티 검색 수";         New = "Chai에 대한 온라인 검색" }
)

foreach ($item in $replacements) {
    $rng = $d.Content
    $found = $rng.Find.Execute($item.Old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = $item.New
        $rng.Font.Bold = 1
    }
}
